$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (not auto-converted to a number),
# matching the source workbook where every cell is stored as inline text,
# and restore the cell's original (default) style afterwards so no
# formatting side effects are introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Price (column D) updates
Set-TextValue "D2"  "228.07"
Set-TextValue "D3"  "22.66"
Set-TextValue "D4"  "5.521"
Set-TextValue "D5"  "0.05551"
Set-TextValue "D6"  "3.417"
Set-TextValue "D7"  "6.493"
Set-TextValue "D8"  "1.142"
Set-TextValue "D9"  "0.7927"
Set-TextValue "D10" "0.1408"
Set-TextValue "D11" "0.07401"
Set-TextValue "D12" "0.03162"
Set-TextValue "D13" "0.02941"
Set-TextValue "D14" "0.09265"
Set-TextValue "D15" "0.001666"
Set-TextValue "D16" "3.272"
Set-TextValue "D17" "0.04730"
Set-TextValue "D18" "0.0005944"
Set-TextValue "E18" "17OneONE"
Set-TextValue "D19" "0.006236"
Set-TextValue "D20" "0.005245"
Set-TextValue "D22" "0.0001508"
Set-TextValue "D23" "3.675"
Set-TextValue "D24" "2.195"
Set-TextValue "D27" "0.0008358"
Set-TextValue "D40" "0.04039"
Set-TextValue "D41" "0.007108"

# Rows 42 & 43 swapped content (BKEXToken <-> CEJI) plus updated prices
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003621"
Set-TextValue "E42" "41CEJICEJI"

Set-TextValue "B43" "BKEXToken"
Set-TextValue "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1028"
Set-TextValue "E43" "42BKEXTokenBKK"

Set-TextValue "D44" "0.008136"
Set-TextValue "E44" "43LocalTradersLCTWorstin24h"

Set-TextValue "D46" "0.00005522"
Set-TextValue "D47" "0.00000000754"
Set-TextValue "D48" "0.6791"
Set-TextValue "D49" "0.09325"
Set-TextValue "D50" "0.00002112"
